$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Frostbolt -> Firebolt, keep "Deal 6(9) Fire damage." ---
$ws.Range("A3").Value = "Firebolt"

# --- Row 4: Shockbolt stays, new description ---
$ws.Range("C4").Value = "Deal 5 Shock damage. Apply 1(2) Vulnerable."

# --- Row 5: Firebolt -> Frostbolt, new description ---
$ws.Range("A5").Value = "Frostbolt"
$ws.Range("C5").Value = "Deal 5(7) Frost damage. Gain 5(7) Block."

# --- Row 15: Fire Rune gets a new, distinct description + Code "x" ---
$ws.Range("C15").Value = "Gain 4(7) Block. The next time you're attacked, deal 4(7) Fire damage back."
$ws.Range("F15").Value = "x"

# --- Row 18: new card "Simon Says" (moved up from the old row 21) ---
$ws.Range("A18").Value = "Simon Says"
$ws.Range("B18").Value = "Power"
$ws.Range("C18").Value = "Each time you deal spell damage, this power chooses Fire, Frost, or Shock randomly. If the next damage you deal is of that type, gain 1 Strength. If not, lose all Strength."
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = "Rare"
$ws.Range("F18").Value = "x"

# --- Row 19: new card "Cold Burn" ---
$ws.Range("A19").Value = "Cold Burn"
$ws.Range("B19").Value = "Attack"
$ws.Range("C19").Value = "Deal 5(7) Frost damage. Deal 5(7) Fire damage."
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = "Uncommon"
$ws.Range("F19").Value = "x"

# --- Row 20: new card "Hypothermia" ---
$ws.Range("A20").Value = "Hypothermia"
$ws.Range("B20").Value = "Skill"
$ws.Range("C20").Value = "Gain 10 Block. Shuffle a Void into your draw pile."
$ws.Range("D20").Value = 0
$ws.Range("E20").Value = "Common"

# --- Row 21: the old "Simon Says" card row is now just the Table2 totals row ---
$ws.Range("A21:E21").ClearContents()

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("G18").Select()

Write-Output "edits applied"
